$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 1-11 (columns B,C,D) with refined values ---
$ws.Range("B1").Value = 3.997212021557795
$ws.Range("C1").Value = 3.3711028133256
$ws.Range("D1").Value = 3.614855494484687

$ws.Range("B2").Value = 3217126.138732513
$ws.Range("C2").Value = 3202355.757565538
$ws.Range("D2").Value = 3208599.328284424

$ws.Range("B3").Value = 804842.5055718543
$ws.Range("C3").Value = 949943.0705307999
$ws.Range("D3").Value = 887614.8253173321

$ws.Range("B4").Value = 0.09665226416387108
$ws.Range("C4").Value = 0.09620851694328061
$ws.Range("D4").Value = 0.09639609281703397

$ws.Range("B5").Value = 0.02417991931441348
$ws.Range("C5").Value = 0.02853918206320463
$ws.Range("D5").Value = 0.02666665181059351

$ws.Range("B6").Value = 2627.936930989841
$ws.Range("C6").Value = 2627.936930989841
$ws.Range("D6").Value = 2627.936930989841

$ws.Range("B7").Value = 10777.64271445416
$ws.Range("C7").Value = 10777.64271445416
$ws.Range("D7").Value = 10777.64271445416

$ws.Range("B8").Value = 0.08000209779143486
$ws.Range("C8").Value = 0.08000209779143486
$ws.Range("D8").Value = 0.08000209779143486

$ws.Range("B9").Value = 33.19535351221312
$ws.Range("C9").Value = 33.19535351221312
$ws.Range("D9").Value = 33.19535351221312

$ws.Range("B10").Value = 79.99994760609384
$ws.Range("C10").Value = 79.99994760609384
$ws.Range("D10").Value = 79.99994760609384

$ws.Range("B11").Value = 0.4149421906581915
$ws.Range("C11").Value = 0.4149421906581915
$ws.Range("D11").Value = 0.4149421906581915

# --- Add new rows 12-15 for base plot properties ---
$ws.Range("A12").Value = "Cp_base"
$ws.Range("B12").Value = -0.04297219340437737
$ws.Range("C12").Value = -0.04297219340437737
$ws.Range("D12").Value = -0.04297219340437737

$ws.Range("A13").Value = "P_base"
$ws.Range("B13").Value = 179.4863027032873
$ws.Range("C13").Value = 179.4863027032873
$ws.Range("D13").Value = 179.4863027032873

$ws.Range("A14").Value = "D_base"
$ws.Range("B14").Value = 37052.9710246601
$ws.Range("C14").Value = 37052.9710246601
$ws.Range("D14").Value = 37052.9710246601

$ws.Range("A15").Value = "S_base"
$ws.Range("B15").Value = 206.4389898649435
$ws.Range("C15").Value = 206.4389898649435
$ws.Range("D15").Value = 206.4389898649435
